# Auto-generated script to apply scheduled market-data refresh
# to the Tonberry_Profits (Leve profit tracker) workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) for the rows whose market data changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 373.5
$ws.Range("I2").Value = 213.5
$ws.Range("J2").Value = 613.5
$ws.Range("K2").Value = 213.5
$ws.Range("L2").Value = 613.5
$ws.Range("M2").Value = -100.5
$ws.Range("N2").Value = -839.5
$ws.Range("H33").Value = 89
$ws.Range("I33").Value = 79.833336
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 79.833336
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 149.166664
$ws.Range("N33").Value = -558
$ws.Range("H40").Value = 2799.8
$ws.Range("J40").Value = 2999.6667
$ws.Range("L40").Value = 2999.6667
$ws.Range("N40").Value = -3349.6667
$ws.Range("H52").Value = 3688.0715
$ws.Range("I52").Value = 3162.6
$ws.Range("K52").Value = 9487.799999999999
$ws.Range("M52").Value = -9327.799999999999
$ws.Range("H62").Value = 2099.5
$ws.Range("I62").Value = 1799.3334
$ws.Range("K62").Value = 1799.3334
$ws.Range("M62").Value = -1175.3334
$ws.Range("H65").Value = 2099.5
$ws.Range("I65").Value = 1799.3334
$ws.Range("K65").Value = 8996.666999999999
$ws.Range("M65").Value = -5876.666999999999
$ws.Range("H104").Value = 1940
$ws.Range("I104").Value = 2292
$ws.Range("K104").Value = 6876
$ws.Range("M104").Value = -5129
$ws.Range("H138").Value = 2302.318
$ws.Range("I138").Value = 2196.8333
$ws.Range("J138").Value = 2528.3572
$ws.Range("K138").Value = 6590.499899999999
$ws.Range("L138").Value = 7585.071599999999
$ws.Range("M138").Value = -1450.499899999999
$ws.Range("N138").Value = -17865.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4492.4565
$ws.Range("I32").Value = 3056.6667
$ws.Range("J32").Value = 8137.154
$ws.Range("K32").Value = 3056.6667
$ws.Range("L32").Value = 8137.154
$ws.Range("M32").Value = -2769.6667
$ws.Range("N32").Value = -8711.154
$ws.Range("H45").Value = 1474.7858
$ws.Range("I45").Value = 909.1667
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 909.1667
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -532.1667
$ws.Range("N45").Value = -2653
$ws.Range("H74").Value = 1225.3914
$ws.Range("J74").Value = 1487.8572
$ws.Range("L74").Value = 1487.8572
$ws.Range("N74").Value = -3235.8572
$ws.Range("H77").Value = 1225.3914
$ws.Range("J77").Value = 1487.8572
$ws.Range("L77").Value = 7439.286
$ws.Range("N77").Value = -16175.286
$ws.Range("H132").Value = 1714.1777
$ws.Range("I132").Value = 1264.7
$ws.Range("J132").Value = 2613.1333
$ws.Range("K132").Value = 3794.1
$ws.Range("L132").Value = 7839.3999
$ws.Range("M132").Value = -1264.1
$ws.Range("N132").Value = -12899.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 269.66666
$ws.Range("I22").Value = 244.5
$ws.Range("K22").Value = 244.5
$ws.Range("M22").Value = -71.5
$ws.Range("H99").Value = 1846.6666
$ws.Range("J99").Value = 2171.4285
$ws.Range("L99").Value = 2171.4285
$ws.Range("N99").Value = -5167.4285
$ws.Range("H105").Value = 2386.3635
$ws.Range("I105").Value = 2386.3635
$ws.Range("K105").Value = 2386.3635
$ws.Range("M105").Value = -639.3634999999999
$ws.Range("H107").Value = 1011
$ws.Range("I107").Value = 1011
$ws.Range("K107").Value = 1011
$ws.Range("M107").Value = 909
$ws.Range("H122").Value = 54000
$ws.Range("J122").Value = 54000
$ws.Range("L122").Value = 54000
$ws.Range("N122").Value = -63800
$ws.Range("H134").Value = 3952.238
$ws.Range("I134").Value = 4157.737
$ws.Range("K134").Value = 12473.211
$ws.Range("M134").Value = -9938.210999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1000.9091
$ws.Range("J22").Value = 1238.75
$ws.Range("L22").Value = 1238.75
$ws.Range("N22").Value = -1938.75
$ws.Range("H31").Value = 1807.8235
$ws.Range("I31").Value = 1358
$ws.Range("K31").Value = 1358
$ws.Range("M31").Value = -1063
$ws.Range("H34").Value = 1807.8235
$ws.Range("I34").Value = 1358
$ws.Range("K34").Value = 1358
$ws.Range("M34").Value = -1156
$ws.Range("H62").Value = 2566
$ws.Range("I62").Value = 2566
$ws.Range("K62").Value = 2566
$ws.Range("M62").Value = -1942
$ws.Range("H65").Value = 2566
$ws.Range("I65").Value = 2566
$ws.Range("K65").Value = 12830
$ws.Range("M65").Value = -9710
$ws.Range("H86").Value = 333335580
$ws.Range("J86").Value = 3399
$ws.Range("L86").Value = 3399
$ws.Range("N86").Value = -5645
$ws.Range("H89").Value = 333335580
$ws.Range("J89").Value = 3399
$ws.Range("L89").Value = 16995
$ws.Range("N89").Value = -28227

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 717.25
$ws.Range("I34").Value = 184.5
$ws.Range("K34").Value = 553.5
$ws.Range("M34").Value = -469.5
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H55").Value = 3999.5
$ws.Range("J55").Value = 3999.5
$ws.Range("L55").Value = 11998.5
$ws.Range("N55").Value = -12352.5
$ws.Range("H122").Value = 1087.7858
$ws.Range("J122").Value = 1131.8182
$ws.Range("L122").Value = 10186.3638
$ws.Range("N122").Value = -15086.3638
$ws.Range("H131").Value = 768.5361
$ws.Range("J131").Value = 789.91113
$ws.Range("L131").Value = 2369.73339
$ws.Range("N131").Value = -12449.73339
$ws.Range("H137").Value = 2477.5454
$ws.Range("I137").Value = 1279.375
$ws.Range("J137").Value = 3162.2144
$ws.Range("K137").Value = 3838.125
$ws.Range("L137").Value = 9486.643199999999
$ws.Range("M137").Value = 1261.875
$ws.Range("N137").Value = -19686.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5700625.5
$ws.Range("I12").Value = 7000000
$ws.Range("J12").Value = 3535001.8
$ws.Range("K12").Value = 7000000
$ws.Range("L12").Value = 3535001.8
$ws.Range("M12").Value = -6999860
$ws.Range("N12").Value = -3535281.8
$ws.Range("H15").Value = 15999.5
$ws.Range("J15").Value = 15999.5
$ws.Range("L15").Value = 15999.5
$ws.Range("N15").Value = -16575.5
$ws.Range("H81").Value = 15999.5
$ws.Range("J81").Value = 15999.5
$ws.Range("L81").Value = 15999.5
$ws.Range("N81").Value = -17995.5
$ws.Range("H84").Value = 15999.5
$ws.Range("J84").Value = 15999.5
$ws.Range("L84").Value = 47998.5
$ws.Range("N84").Value = -57982.5
$ws.Range("H132").Value = 2962772
$ws.Range("I132").Value = 4276891.5
$ws.Range("J132").Value = 6002.75
$ws.Range("K132").Value = 12830674.5
$ws.Range("L132").Value = 18008.25
$ws.Range("M132").Value = -12828144.5
$ws.Range("N132").Value = -23068.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3371.7222
$ws.Range("I7").Value = 2142.6667
$ws.Range("K7").Value = 2142.6667
$ws.Range("M7").Value = -2030.6667
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H20").Value = 9185.454
$ws.Range("I20").Value = 7808
$ws.Range("J20").Value = 10333.333
$ws.Range("K20").Value = 7808
$ws.Range("L20").Value = 10333.333
$ws.Range("M20").Value = -7582
$ws.Range("N20").Value = -10785.333
$ws.Range("H55").Value = 289.3793
$ws.Range("I55").Value = 210.08696
$ws.Range("J55").Value = 593.3333
$ws.Range("K55").Value = 210.08696
$ws.Range("L55").Value = 593.3333
$ws.Range("M55").Value = -37.08696
$ws.Range("N55").Value = -939.3333
$ws.Range("H82").Value = 2026.2
$ws.Range("I82").Value = 1182.7142
$ws.Range("J82").Value = 3994.3333
$ws.Range("K82").Value = 1182.7142
$ws.Range("L82").Value = 3994.3333
$ws.Range("M82").Value = -821.7141999999999
$ws.Range("N82").Value = -4716.3333
$ws.Range("H85").Value = 2026.2
$ws.Range("I85").Value = 1182.7142
$ws.Range("J85").Value = 3994.3333
$ws.Range("K85").Value = 1182.7142
$ws.Range("L85").Value = 3994.3333
$ws.Range("M85").Value = 65.28580000000011
$ws.Range("N85").Value = -6490.3333
$ws.Range("H126").Value = 3371.7222
$ws.Range("I126").Value = 2142.6667
$ws.Range("K126").Value = 6428.000100000001
$ws.Range("M126").Value = -3958.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 80005
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H126").Value = 2811.8333
$ws.Range("I126").Value = 1419.4615
$ws.Range("J126").Value = 6432
$ws.Range("K126").Value = 4258.3845
$ws.Range("L126").Value = 19296
$ws.Range("M126").Value = -1788.3845
$ws.Range("N126").Value = -24236
$ws.Range("H132").Value = 1871.75
$ws.Range("I132").Value = 1279.7693
$ws.Range("K132").Value = 3839.3079
$ws.Range("M132").Value = -1309.3079
